$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "AOF3JC"
$ws.Range("B26").Value = "Cinta Flex de cabezal para Epson"
$ws.Range("C26").Value = "R260 R360 R380 R390 RX580 RX590 R1390 R1400 R1410 R1430 1500w L1800 R1800 R1900 R2000 R2400 R2880 EP4004"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 100000
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 0
$ws.Range("H26").Formula = "=(E26-D26)*G26"
$ws.Range("I26").Formula = "=D26*F26"
$ws.Range("J26").Value = 0
